$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 129 and 130: fill in resultado / profit ---
$ws.Cells.Item(129, 7).Value = "Acierto"
$ws.Cells.Item(129, 8).Value = 0.53
$ws.Cells.Item(130, 7).Value = "Acierto"
$ws.Cells.Item(130, 8).Value = 0.53

# --- Append new rows 131-136 ---
# Use row 130 (which already has the "2025-10-10" text date and empty
# resultado/profit placeholders) as a paste-values template so the new
# date cells land as text instead of being auto-converted to date serials.
$ws.Range("A130:H130").Copy()
for ($r = 131; $r -le 136; $r++) {
    $ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4163)
}
$excel.CutCopyMode = 0

$ws.Cells.Item(131, 1).Value = 14807180
$ws.Cells.Item(131, 3).Value = "Andre Ilagan"
$ws.Cells.Item(131, 4).Value = "Dominique Rolland"
$ws.Cells.Item(131, 5).Value = "Gana Dominique Rolland"
$ws.Cells.Item(131, 6).Value = 3

$ws.Cells.Item(132, 1).Value = 14807184
$ws.Cells.Item(132, 3).Value = "Edward Winter"
$ws.Cells.Item(132, 4).Value = "Abdullah Shelbayh"
$ws.Cells.Item(132, 5).Value = "Gana Edward Winter"
$ws.Cells.Item(132, 6).Value = 2.25

$ws.Cells.Item(133, 1).Value = 14859921
$ws.Cells.Item(133, 3).Value = "Abel Forger"
$ws.Cells.Item(133, 4).Value = "Sean Cuenin"
$ws.Cells.Item(133, 5).Value = "Gana Abel Forger"
$ws.Cells.Item(133, 6).Value = 1.67

$ws.Cells.Item(134, 1).Value = 14859635
$ws.Cells.Item(134, 3).Value = "Robin Bertrand"
$ws.Cells.Item(134, 4).Value = "Emilien Demanet"
$ws.Cells.Item(134, 5).Value = "Gana Emilien Demanet"
$ws.Cells.Item(134, 6).Value = 2.62

$ws.Cells.Item(135, 1).Value = 14859636
$ws.Cells.Item(135, 3).Value = "Branko Djuric"
$ws.Cells.Item(135, 4).Value = "Igor Kudriashov"
$ws.Cells.Item(135, 5).Value = "Gana Igor Kudriashov"
$ws.Cells.Item(135, 6).Value = 1.91

$ws.Cells.Item(136, 1).Value = 14859949
$ws.Cells.Item(136, 3).Value = "Tim Handel"
$ws.Cells.Item(136, 4).Value = "Daniel Michalski"
$ws.Cells.Item(136, 5).Value = "Gana Tim Handel"
$ws.Cells.Item(136, 6).Value = 3.4

# resultado / profit for the new rows stay blank (no result yet), matching
# the source data - clear any inherited template values just in case.
$ws.Range("G131:H136").ClearContents()
